# Updates res_bus/vm_pu.xlsx values for the 380 kV case (Case_2_43).
# Slack-bus voltage setpoint (col B) drops from 1.05 to 1.02 p.u. for all
# timesteps (rows 2-25), and the resulting per-bus voltage magnitudes in
# columns C:F and I:N are updated to the recomputed power-flow results.
# Column G (slack bus itself) stays at 1, column H has no data (gap in the
# network numbering), and column A (timestep index) is unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$r2a = New-Object "object[,]" 1,5
$r2a[0,0] = 1.02
$r2a[0,1] = 1.015497782699737
$r2a[0,2] = 1.022445918058817
$r2a[0,3] = 0.9926147277508489
$r2a[0,4] = 1.027621456089077
$ws.Range("B2:F2").Value = $r2a

$r2b = New-Object "object[,]" 1,6
$r2b[0,0] = 1.028194577939249
$r2b[0,1] = 1.020722958502823
$r2b[0,2] = 1.025280276155554
$r2b[0,3] = 0.9955398523336033
$r2b[0,4] = 1.030440650512469
$r2b[0,5] = 1.011043303306251
$ws.Range("I2:N2").Value = $r2b

$r3a = New-Object "object[,]" 1,5
$r3a[0,0] = 1.02
$r3a[0,1] = 1.016340864875646
$r3a[0,2] = 1.02307307580092
$r3a[0,3] = 0.9936372048519304
$r3a[0,4] = 1.028719273825873
$ws.Range("B3:F3").Value = $r3a

$r3b = New-Object "object[,]" 1,6
$r3b[0,0] = 1.028352070087204
$r3b[0,1] = 1.021201642612424
$r3b[0,2] = 1.025714592281431
$r3b[0,3] = 0.9963617723202692
$r3b[0,4] = 1.031345448279973
$r3b[0,5] = 1.011200814450409
$ws.Range("I3:N3").Value = $r3b

$r4a = New-Object "object[,]" 1,5
$r4a[0,0] = 1.02
$r4a[0,1] = 1.016886416656375
$r4a[0,2] = 1.023478465166813
$r4a[0,3] = 0.9942998659930995
$r4a[0,4] = 1.02942965303633
$ws.Range("B4:F4").Value = $r4a

$r4b = New-Object "object[,]" 1,6
$r4b[0,0] = 1.028452069755527
$r4b[0,1] = 1.021510783991011
$r4b[0,2] = 1.025994513501666
$r4b[0,3] = 0.9968940712668345
$r4b[0,4] = 1.031930347817456
$r4b[0,5] = 1.011302529479239
$ws.Range("I4:N4").Value = $r4b

$r5a = New-Object "object[,]" 1,5
$r5a[0,0] = 1.02
$r5a[0,1] = 1.0171157706328
$r5a[0,2] = 1.023648787953266
$r5a[0,3] = 0.9945786998346017
$r5a[0,4] = 1.029728300296114
$ws.Range("B5:F5").Value = $r5a

$r5b = New-Object "object[,]" 1,6
$r5b[0,0] = 1.028493652057705
$r5b[0,1] = 1.021640602673973
$r5b[0,2] = 1.026111925478718
$r5b[0,3] = 0.997117960005301
$r5b[0,4] = 1.032176103765295
$r5b[0,5] = 1.011345241040793
$ws.Range("I5:N5").Value = $r5b

$r6a = New-Object "object[,]" 1,5
$r6a[0,0] = 1.02
$r6a[0,1] = 1.017154280378291
$r6a[0,2] = 1.023677379820905
$r6a[0,3] = 0.9946255319796338
$r6a[0,4] = 1.029778444698482
$ws.Range("B6:F6").Value = $r6a

$r6b = New-Object "object[,]" 1,6
$r6b[0,0] = 1.028500607056953
$r6b[0,1] = 1.02166239131114
$r6b[0,2] = 1.026131623783283
$r6b[0,3] = 0.9971555583673453
$r6b[0,4] = 1.032217359290478
$r6b[0,5] = 1.01135240959194
$ws.Range("I6:N6").Value = $r6b

$r7a = New-Object "object[,]" 1,5
$r7a[0,0] = 1.02
$r7a[0,1] = 1.016889481281043
$r7a[0,2] = 1.023480741434446
$r7a[0,3] = 0.9943035907982488
$r7a[0,4] = 1.029433643561067
$ws.Range("B7:F7").Value = $r7a

$r7b = New-Object "object[,]" 1,6
$r7b[0,0] = 1.02845262717908
$r7b[0,1] = 1.021512519202944
$r7b[0,2] = 1.025996083416376
$r7b[0,3] = 0.9968970624462087
$r7b[0,4] = 1.031933632154811
$r7b[0,5] = 1.011303100387857
$ws.Range("I7:N7").Value = $r7b

$r8a = New-Object "object[,]" 1,5
$r8a[0,0] = 1.02
$r8a[0,1] = 1.015782701225326
$r8a[0,2] = 1.022657955779921
$r8a[0,3] = 0.9929600610674301
$r8a[0,4] = 1.027992464822022
$ws.Range("B8:F8").Value = $r8a

$r8b = New-Object "object[,]" 1,6
$r8b[0,0] = 1.028248197844678
$r8b[0,1] = 1.020884855583494
$r8b[0,2] = 1.025427284581886
$r8b[0,3] = 0.995817528259106
$r8b[0,4] = 1.030746547851494
$r8b[0,5] = 1.011096577241308
$ws.Range("I8:N8").Value = $r8b

$r9a = New-Object "object[,]" 1,5
$r9a[0,0] = 1.02
$r9a[0,1] = 1.013832623114584
$r9a[0,2] = 1.021204922859877
$r9a[0,3] = 0.9906006454969559
$r9a[0,4] = 1.025453076541555
$ws.Range("B9:F9").Value = $r9a

$r9b = New-Object "object[,]" 1,6
$r9b[0,0] = 1.027873383106093
$r9b[0,1] = 1.019774281502987
$r9b[0,2] = 1.024416531672746
$r9b[0,3] = 0.9939188001724441
$r9b[0,4] = 1.028650448007413
$r9b[0,5] = 1.010731098839352
$ws.Range("I9:N9").Value = $r9b

$r10a = New-Object "object[,]" 1,5
$r10a[0,0] = 1.02
$r10a[0,1] = 1.012532773823903
$r10a[0,2] = 1.020234182992269
$r10a[0,3] = 0.989033133672735
$r10a[0,4] = 1.023760273293209
$ws.Range("B10:F10").Value = $r10a

$r10b = New-Object "object[,]" 1,6
$r10b[0,0] = 1.027613741968303
$r10b[0,1] = 1.019030890536559
$r10b[0,2] = 1.023737073150393
$r10b[0,3] = 0.9926553831429383
$r10b[0,4] = 1.027250180636101
$r10b[0,5] = 1.010486416670771
$ws.Range("I10:N10").Value = $r10b

$r11a = New-Object "object[,]" 1,5
$r11a[0,0] = 1.02
$r11a[0,1] = 1.011969985296794
$r11a[0,2] = 1.019813373237379
$r11a[0,3] = 0.988355674866747
$r11a[0,4] = 1.023027302822807
$ws.Range("B11:F11").Value = $r11a

$r11b = New-Object "object[,]" 1,6
$r11b[0,0] = 1.027499006763025
$r11b[0,1] = 1.018708290970139
$r11b[0,2] = 1.023441539153863
$r11b[0,3] = 0.9921088820399291
$r11b[0,4] = 1.026643174101625
$r11b[0,5] = 1.010380225937092
$ws.Range("I11:N11").Value = $r11b

$r12a = New-Object "object[,]" 1,5
$r12a[0,0] = 1.02
$r12a[0,1] = 1.011760949752831
$r12a[0,2] = 1.019656996095024
$r12a[0,3] = 0.9881042295826724
$r12a[0,4] = 1.022755048753188
$ws.Range("B12:F12").Value = $r12a

$r12b = New-Object "object[,]" 1,6
$r12b[0,0] = 1.02745604273482
$r12b[0,1] = 1.018588357798084
$r12b[0,2] = 1.023331566738024
$r12b[0,3] = 0.9919059725120875
$r12b[0,4] = 1.026417602730466
$r12b[0,5] = 1.010340745928849
$ws.Range("I12:N12").Value = $r12b

$r13a = New-Object "object[,]" 1,5
$r13a[0,0] = 1.02
$r13a[0,1] = 1.01180578815643
$r13a[0,2] = 1.019690542652949
$r13a[0,3] = 0.9881581567098651
$r13a[0,4] = 1.022813448001591
$ws.Range("B13:F13").Value = $r13a

$r13b = New-Object "object[,]" 1,6
$r13b[0,0] = 1.027465274324036
$r13b[0,1] = 1.018614088624283
$r13b[0,2] = 1.023355165135656
$r13b[0,3] = 0.9919494934313052
$r13b[0,4] = 1.026465993171983
$r13b[0,5] = 1.010349216151167
$ws.Range("I13:N13").Value = $r13b

$r14a = New-Object "object[,]" 1,5
$r14a[0,0] = 1.02
$r14a[0,1] = 1.011952706160964
$r14a[0,2] = 1.019800448467064
$r14a[0,3] = 0.9883348863814464
$r14a[0,4] = 1.023004798138334
$ws.Range("B14:F14").Value = $r14a

$r14b = New-Object "object[,]" 1,6
$r14b[0,0] = 1.027495462403182
$r14b[0,1] = 1.018698379400833
$r14b[0,2] = 1.02343245283175
$r14b[0,3] = 0.9920921077337197
$r14b[0,4] = 1.026624530383316
$r14b[0,5] = 1.010376963240938
$ws.Range("I14:N14").Value = $r14b

$r15a = New-Object "object[,]" 1,5
$r15a[0,0] = 1.02
$r15a[0,1] = 1.012043228436076
$r15a[0,2] = 1.019868155877123
$r15a[0,3] = 0.9884438009545853
$r15a[0,4] = 1.023122695779696
$ws.Range("B15:F15").Value = $r15a

$r15b = New-Object "object[,]" 1,6
$r15b[0,0] = 1.027514016410092
$r15b[0,1] = 1.018750299796735
$r15b[0,2] = 1.023480046133569
$r15b[0,3] = 0.9921799884222134
$r15b[0,4] = 1.026722196868172
$r15b[0,5] = 1.010394054371106
$ws.Range("I15:N15").Value = $r15b

$r16a = New-Object "object[,]" 1,5
$r16a[0,0] = 1.02
$r16a[0,1] = 1.012570125479808
$r16a[0,2] = 1.020262100904096
$r16a[0,3] = 0.9890781214508737
$r16a[0,4] = 1.023808918655754
$ws.Range("B16:F16").Value = $r16a

$r16b = New-Object "object[,]" 1,6
$r16b[0,0] = 1.027621307973512
$r16b[0,1] = 1.019052285592879
$r16b[0,2] = 1.023756658957774
$r16b[0,3] = 0.9926916645766087
$r16b[0,4] = 1.027290451325264
$r16b[0,5] = 1.01049345913055
$ws.Range("I16:N16").Value = $r16b

$r17a = New-Object "object[,]" 1,5
$r17a[0,0] = 1.02
$r17a[0,1] = 1.012900649410468
$r17a[0,2] = 1.02050908657011
$r17a[0,3] = 0.989476357848556
$r17a[0,4] = 1.024239374781541
$ws.Range("B17:F17").Value = $r17a

$r17b = New-Object "object[,]" 1,6
$r17b[0,0] = 1.027687991422724
$r17b[0,1] = 1.019241524751769
$r17b[0,2] = 1.023929817104597
$r17b[0,3] = 0.9930127773699352
$r17b[0,4] = 1.027646720065462
$r17b[0,5] = 1.010555748585527
$ws.Range("I17:N17").Value = $r17b

$r18a = New-Object "object[,]" 1,5
$r18a[0,0] = 1.02
$r18a[0,1] = 1.013093443651131
$r18a[0,2] = 1.020653103400962
$r18a[0,3] = 0.9897087662937556
$r18a[0,4] = 1.024490454913995
$ws.Range("B18:F18").Value = $r18a

$r18b = New-Object "object[,]" 1,6
$r18b[0,0] = 1.027726663935935
$r18b[0,1] = 1.019351836543815
$r18b[0,2] = 1.024030689512565
$r18b[0,3] = 0.9932001317071769
$r18b[0,4] = 1.027854459818964
$r18b[0,5] = 1.01059205762789
$ws.Range("I18:N18").Value = $r18b

$r19a = New-Object "object[,]" 1,5
$r19a[0,0] = 1.02
$r19a[0,1] = 1.013159182362316
$r19a[0,2] = 1.02070220161935
$r19a[0,3] = 0.9897880325774034
$r19a[0,4] = 1.024576067146836
$ws.Range("B19:F19").Value = $r19a

$r19b = New-Object "object[,]" 1,6
$r19b[0,0] = 1.027739812451436
$r19b[0,1] = 1.019389438420706
$r19b[0,2] = 1.024065062688228
$r19b[0,3] = 0.9932640239640975
$r19b[0,4] = 1.02792528254143
$r19b[0,5] = 1.010604434105208
$ws.Range("I19:N19").Value = $r19b

$r20a = New-Object "object[,]" 1,5
$r20a[0,0] = 1.02
$r20a[0,1] = 1.012865186773553
$r20a[0,2] = 1.020482592065454
$r20a[0,3] = 0.9894336180360679
$r20a[0,4] = 1.024193190668243
$ws.Range("B20:F20").Value = $r20a

$r20b = New-Object "object[,]" 1,6
$r20b[0,0] = 1.027680859961949
$r20b[0,1] = 1.019221228230205
$r20b[0,2] = 1.023911252086986
$r20b[0,3] = 0.9929783193494215
$r20b[0,4] = 1.027608502602641
$r20b[0,5] = 1.010549067928953
$ws.Range("I20:N20").Value = $r20b

$r21a = New-Object "object[,]" 1,5
$r21a[0,0] = 1.02
$r21a[0,1] = 1.011909442209703
$r21a[0,2] = 1.019768085855998
$r21a[0,3] = 0.9882828385668249
$r21a[0,4] = 1.022948450186566
$ws.Range("B21:F21").Value = $r21a

$r21b = New-Object "object[,]" 1,6
$r21b[0,0] = 1.027486582319944
$r21b[0,1] = 1.018673560769281
$r21b[0,2] = 1.023409698987619
$r21b[0,3] = 0.9920501090198102
$r21b[0,4] = 1.02657784794278
$r21b[0,5] = 1.010368793407288
$ws.Range("I21:N21").Value = $r21b

$r22a = New-Object "object[,]" 1,5
$r22a[0,0] = 1.02
$r22a[0,1] = 1.011308581205029
$r22a[0,2] = 1.01931844508019
$r22a[0,3] = 0.9875604150241495
$r22a[0,4] = 1.022165854474109
$ws.Range("B22:F22").Value = $r22a

$r22b = New-Object "object[,]" 1,6
$r22b[0,0] = 1.027362429178737
$r22b[0,1] = 1.018328611486051
$r22b[0,2] = 1.023093207960133
$r22b[0,3] = 0.9914670000341481
$r22b[0,4] = 1.025929243764113
$r22b[0,5] = 1.010255239270498
$ws.Range("I22:N22").Value = $r22b

$r23a = New-Object "object[,]" 1,5
$r23a[0,0] = 1.02
$r23a[0,1] = 1.011627103604551
$r23a[0,2] = 1.019556845838294
$r23a[0,3] = 0.9879432794643023
$r23a[0,4] = 1.022580721142451
$ws.Range("B23:F23").Value = $r23a

$r23b = New-Object "object[,]" 1,6
$r23b[0,0] = 1.027428434771446
$r23b[0,1] = 1.018511533091263
$r23b[0,2] = 1.023261094019752
$r23b[0,3] = 0.991776070289318
$r23b[0,4] = 1.026273137011645
$r23b[0,5] = 1.010315456127898
$ws.Range("I23:N23").Value = $r23b

$r24a = New-Object "object[,]" 1,5
$r24a[0,0] = 1.02
$r24a[0,1] = 1.012881210796605
$r24a[0,2] = 1.020494563933087
$r24a[0,3] = 0.9894529299347244
$r24a[0,4] = 1.024214059273627
$ws.Range("B24:F24").Value = $r24a

$r24b = New-Object "object[,]" 1,6
$r24b[0,0] = 1.027684083050595
$r24b[0,1] = 1.019230399565445
$r24b[0,2] = 1.023919641214301
$r24b[0,3] = 0.9929938892766442
$r24b[0,4] = 1.027625771632152
$r24b[0,5] = 1.010552086702237
$ws.Range("I24:N24").Value = $r24b

$r25a = New-Object "object[,]" 1,5
$r25a[0,0] = 1.02
$r25a[0,1] = 1.014336733809736
$r25a[0,2] = 1.021580934533696
$r25a[0,3] = 0.9912096547607049
$r25a[0,4] = 1.026109549137196
$ws.Range("B25:F25").Value = $r25a

$r25b = New-Object "object[,]" 1,6
$r25b[0,0] = 1.027972005756369
$r25b[0,1] = 1.020061925479589
$r25b[0,2] = 1.024678831209658
$r25b[0,3] = 0.9944092447426414
$r25b[0,4] = 1.029192847452159
$r25b[0,5] = 1.01082576658615
$ws.Range("I25:N25").Value = $r25b

